$wb = $excel.ActiveWorkbook

# --- create the two new sheets, positioned after tool_bimanual ---
$wsTaskTime = $wb.Worksheets.Item("task_time")
$wsBimanual = $wb.Worksheets.Item("tool_bimanual")
$wsPupilDilation = $wb.Worksheets.Add($null, $wsBimanual)
$wsPupilDilation.Name = "pupil_dilation"
$wsPupilBlinks = $wb.Worksheets.Add($null, $wsPupilDilation)
$wsPupilBlinks.Name = "pupil_blinks"

# --- pupil_dilation: header row ---
$wsPupilDilation.Range("A1").Value = "i"
$wsPupilDilation.Range("B1").Value = "Author"
$wsPupilDilation.Range("C1").Value = "Year"
$wsPupilDilation.Range("D1").Value = "Study"
$wsPupilDilation.Range("E1").Value = "Journal"
$wsPupilDilation.Range("F1").Value = "technique"
$wsPupilDilation.Range("G1").Value = "task"
$wsPupilDilation.Range("H1").Value = "note"
$wsPupilDilation.Range("I1").Value = "Nn"
$wsPupilDilation.Range("J1").Value = "Mn"
$wsPupilDilation.Range("K1").Value = "SDn"
$wsPupilDilation.Range("L1").Value = "Ne"
$wsPupilDilation.Range("M1").Value = "Me"
$wsPupilDilation.Range("N1").Value = "SDe"
$wsPupilDilation.Range("O1").Value = "SDpooled"
$wsPupilDilation.Range("P1").Value = "SMD"
$wsPupilDilation.Range("Q1").Value = "g"
$wsPupilDilation.Range("R1").Value = "SDg"

# --- pupil_dilation row 2: Castner et al. ---
$wsPupilDilation.Range("A2").Value = 0
$wsPupilDilation.Range("B2").Value = "Castner et al."
$wsPupilDilation.Range("C2").Value = 2020
$wsPupilDilation.Range("D2").Value = "Pupil diameter differentiates expertise in dental radiography visual search"
$wsPupilDilation.Range("E2").Value = "PLOS ONE"
$wsPupilDilation.Range("F2").Value = "Radiography"
$wsPupilDilation.Range("G2").Value = "Dental radiography, visual search"
$wsPupilDilation.Range("H2").Value = "Reported values are medians? Median change from baseline"
$wsPupilDilation.Range("I2").Formula = "=50*20"
$wsPupilDilation.Range("J2").Value = 0.314
$wsPupilDilation.Range("K2").Value = 0.315
$wsPupilDilation.Range("L2").Formula = "=26*15"
$wsPupilDilation.Range("M2").Value = 0.057
$wsPupilDilation.Range("N2").Value = 0.353
$wsPupilDilation.Range("O2").Formula = "=SQRT(((I2-1)*POWER(K2,2) + (L2-1)*POWER(N2,2))/((I2-1)+(L2-1)))"
$wsPupilDilation.Range("P2").Formula = "=(J2-M2)/O2"
$wsPupilDilation.Range("Q2").Formula = "=P2*(1- (3/(4*(I2+L2)-9)))"
$wsPupilDilation.Range("R2").Formula = "=SQRT((I2+L2)/(I2*L2)+(POWER(P2,2)/(2*(I2+L2))))"

# --- pupil_dilation row 3: Cabrera-Mino et al. ---
$wsPupilDilation.Range("A3").Value = 1
$wsPupilDilation.Range("B3").Value = "Cabrera-Mino et al."
$wsPupilDilation.Range("C3").Value = 2019
$wsPupilDilation.Range("D3").Value = "Task-Evoked Pupillary Responses in Nursing Simulation as an Indicator of Stress and Cognitive Load"
$wsPupilDilation.Range("E3").Value = "Clinical Simulation in Nursing"
$wsPupilDilation.Range("F3").Value = "Various nursing tasks"
$wsPupilDilation.Range("G3").Value = "Elevate HOB"
$wsPupilDilation.Range("H3").Value = "There were different tasks, picked the one that had the most significant result. Values estimated from barplot"
$wsPupilDilation.Range("I3").Value = 13
$wsPupilDilation.Range("J3").Value = 0.75
$wsPupilDilation.Range("K3").Value = 0.75
$wsPupilDilation.Range("L3").Value = 15
$wsPupilDilation.Range("M3").Value = 0.25
$wsPupilDilation.Range("N3").Value = 0.4
$wsPupilDilation.Range("O3").Formula = "=SQRT(((I3-1)*POWER(K3,2) + (L3-1)*POWER(N3,2))/((I3-1)+(L3-1)))"
$wsPupilDilation.Range("P3").Formula = "=(J3-M3)/O3"
$wsPupilDilation.Range("Q3").Formula = "=P3*(1- (3/(4*(I3+L3)-9)))"
$wsPupilDilation.Range("R3").Formula = "=SQRT((I3+L3)/(I3*L3)+(POWER(P3,2)/(2*(I3+L3))))"

# --- pupil_dilation row 4: Bednarik et al. ---
$wsPupilDilation.Range("A4").Value = 2
$wsPupilDilation.Range("B4").Value = "Bednarik et al."
$wsPupilDilation.Range("C4").Value = 2018
$wsPupilDilation.Range("D4").Value = "Pupil Size As an Indicator of Visual-motor Workload and Expertise in Microsurgical Training Tasks"
$wsPupilDilation.Range("E4").Value = "Proceedings of the 2018 ACM Symposium on Eye Tracking Research & Applications"
$wsPupilDilation.Range("F4").Value = "Microsurgery"
$wsPupilDilation.Range("G4").Value = "Suturing"
$wsPupilDilation.Range("H4").Value = "Took the segment 'needle push', estimated from plots"
$wsPupilDilation.Range("I4").Value = 50
$wsPupilDilation.Range("J4").Formula = "=0.02/8"
$wsPupilDilation.Range("K4").Value = 0.005
$wsPupilDilation.Range("L4").Value = 60
$wsPupilDilation.Range("M4").Value = 0.0175
$wsPupilDilation.Range("N4").Value = 0.005
$wsPupilDilation.Range("O4").Formula = "=SQRT(((I4-1)*POWER(K4,2) + (L4-1)*POWER(N4,2))/((I4-1)+(L4-1)))"
$wsPupilDilation.Range("P4").Formula = "=(J4-M4)/O4"
$wsPupilDilation.Range("Q4").Formula = "=P4*(1- (3/(4*(I4+L4)-9)))"
$wsPupilDilation.Range("R4").Formula = "=SQRT((I4+L4)/(I4*L4)+(POWER(P4,2)/(2*(I4+L4))))"

# --- pupil_dilation row 5: Gunawardena et al. ---
$wsPupilDilation.Range("A5").Value = 3
$wsPupilDilation.Range("B5").Value = "Gunawardena et al."
$wsPupilDilation.Range("C5").Value = 2019
$wsPupilDilation.Range("D5").Value = "Assessing Surgeons’ Skill Level in Laparoscopic Cholecystectomy using Eye Metrics"
$wsPupilDilation.Range("E5").Value = "Eye Tracking Research and Applications Symposium (ETRA)"
$wsPupilDilation.Range("F5").Value = "Laparoscopy"
$wsPupilDilation.Range("G5").Value = "Laparoscopic cholecystectomy"
$wsPupilDilation.Range("H5").Value = "Study had only 4 participants of 3 skill levels who completed >=7 tasks each. I picked the least experienced participant and expert E-2."
$wsPupilDilation.Range("I5").Value = 7
$wsPupilDilation.Range("J5").Value = 4.87
$wsPupilDilation.Range("K5").Value = 0.56
$wsPupilDilation.Range("L5").Value = 7
$wsPupilDilation.Range("M5").Value = 4.1
$wsPupilDilation.Range("N5").Value = 0.31
$wsPupilDilation.Range("O5").Formula = "=SQRT(((I5-1)*POWER(K5,2) + (L5-1)*POWER(N5,2))/((I5-1)+(L5-1)))"
$wsPupilDilation.Range("P5").Formula = "=(J5-M5)/O5"
$wsPupilDilation.Range("Q5").Formula = "=P5*(1- (3/(4*(I5+L5)-9)))"
$wsPupilDilation.Range("R5").Formula = "=SQRT((I5+L5)/(I5*L5)+(POWER(P5,2)/(2*(I5+L5))))"

# --- task_time: add row 16 (Zheng et al.) ---
$wsTaskTime.Range("A16").Value = 14
$wsTaskTime.Range("B16").Value = "Zheng et al."
$wsTaskTime.Range("C16").Value = 2021
$wsTaskTime.Range("D16").Value = "Action-related eye measures to assess surgical expertise"
$wsTaskTime.Range("E16").Value = "BJS Open"
$wsTaskTime.Range("F16").Value = "Laparoscopy"
$wsTaskTime.Range("G16").Value = "Box trainer"
$wsTaskTime.Range("H16").Value = "Transporting and loading task"
$wsTaskTime.Range("I16").Value = 12
$wsTaskTime.Range("J16").Value = 6.296
$wsTaskTime.Range("K16").Value = 1.853
$wsTaskTime.Range("L16").Value = 5
$wsTaskTime.Range("M16").Value = 2.96
$wsTaskTime.Range("N16").Value = 0.752
$wsTaskTime.Range("O16").Formula = "=SQRT(((I16-1)*POWER(K16,2) + (L16-1)*POWER(N16,2))/((I16-1)+(L16-1)))"
$wsTaskTime.Range("P16").Formula = "=(J16-M16)/O16"
$wsTaskTime.Range("Q16").Formula = "=P16*(1- (3/(4*(I16+L16)-9)))"
$wsTaskTime.Range("R16").Formula = "=SQRT((I16+L16)/(I16*L16)+(POWER(P16,2)/(2*(I16+L16))))"

# --- pupil_dilation row 6: Dilley et al. ---
$wsPupilDilation.Range("A6").Value = 4
$wsPupilDilation.Range("D6").Value = "Visual behaviour in robotic surgery—Demonstrating the validity of the simulated environment"
$wsPupilDilation.Range("B6").Value = "Dilley et al."
$wsPupilDilation.Range("C6").Value = 2020
$wsPupilDilation.Range("E6").Value = "International Journal of Medical Robotics and Computer Assisted Surgery"
$wsPupilDilation.Range("F6").Value = "Robotic surgery"
$wsPupilDilation.Range("G6").Value = "Fundamentals of Robotic Surgery, simulator task"
$wsPupilDilation.Range("H6").Value = "SDs calculated from inter-quartile ranges (SD = (3/4)*IQR). The paper reports medians."
$wsPupilDilation.Range("I6").Value = 18
$wsPupilDilation.Range("J6").Value = 3.25
$wsPupilDilation.Range("K6").Formula = "=0.96*(3/4)"
$wsPupilDilation.Range("L6").Value = 14
$wsPupilDilation.Range("M6").Value = 3.26
$wsPupilDilation.Range("N6").Formula = "=0.7*(3/4)"
$wsPupilDilation.Range("O6").Formula = "=SQRT(((I6-1)*POWER(K6,2) + (L6-1)*POWER(N6,2))/((I6-1)+(L6-1)))"
$wsPupilDilation.Range("P6").Formula = "=(J6-M6)/O6"
$wsPupilDilation.Range("Q6").Formula = "=P6*(1- (3/(4*(I6+L6)-9)))"
$wsPupilDilation.Range("R6").Formula = "=SQRT((I6+L6)/(I6*L6)+(POWER(P6,2)/(2*(I6+L6))))"

# --- pupil_dilation row 7: Gao et al. ---
$wsPupilDilation.Range("A7").Value = 5
$wsPupilDilation.Range("B7").Value = "Gao et al."
$wsPupilDilation.Range("C7").Value = 2018
$wsPupilDilation.Range("E7").Value = "American Surgeon"
$wsPupilDilation.Range("D7").Value = "Quantitative evaluations of the effects of noise on mental workloads based on pupil dilation during laparoscopic surgery"
$wsPupilDilation.Range("F7").Value = "Laparoscopy"
$wsPupilDilation.Range("G7").Value = "Appendectromy simulator"
$wsPupilDilation.Range("H7").Value = "They evaluated different noise conditions, I picked values from the no-noise condition. Paper does not give explicitly the number of participants in groups, only total number (24) which was `"divided into experienced and moderately experienced`". I assumed 12 per group"
$wsPupilDilation.Range("I7").Value = 12
$wsPupilDilation.Range("J7").Value = 0.108
$wsPupilDilation.Range("K7").Value = 0.075
$wsPupilDilation.Range("L7").Value = 12
$wsPupilDilation.Range("M7").Value = 0.038
$wsPupilDilation.Range("N7").Value = 0.023
$wsPupilDilation.Range("O7").Formula = "=SQRT(((I7-1)*POWER(K7,2) + (L7-1)*POWER(N7,2))/((I7-1)+(L7-1)))"
$wsPupilDilation.Range("P7").Formula = "=(J7-M7)/O7"
$wsPupilDilation.Range("Q7").Formula = "=P7*(1- (3/(4*(I7+L7)-9)))"
$wsPupilDilation.Range("R7").Formula = "=SQRT((I7+L7)/(I7*L7)+(POWER(P7,2)/(2*(I7+L7))))"

# --- pupil_blinks: header row ---
$wsPupilBlinks.Range("A1").Value = "i"
$wsPupilBlinks.Range("B1").Value = "Author"
$wsPupilBlinks.Range("C1").Value = "Year"
$wsPupilBlinks.Range("D1").Value = "Study"
$wsPupilBlinks.Range("E1").Value = "Journal"
$wsPupilBlinks.Range("F1").Value = "technique"
$wsPupilBlinks.Range("G1").Value = "task"
$wsPupilBlinks.Range("H1").Value = "note"
$wsPupilBlinks.Range("I1").Value = "Nn"
$wsPupilBlinks.Range("J1").Value = "Mn"
$wsPupilBlinks.Range("K1").Value = "SDn"
$wsPupilBlinks.Range("L1").Value = "Ne"
$wsPupilBlinks.Range("M1").Value = "Me"
$wsPupilBlinks.Range("N1").Value = "SDe"
$wsPupilBlinks.Range("O1").Value = "SDpooled"
$wsPupilBlinks.Range("P1").Value = "SMD"
$wsPupilBlinks.Range("Q1").Value = "g"
$wsPupilBlinks.Range("R1").Value = "SDg"

# --- pupil_blinks row 2: Dilley et al. ---
$wsPupilBlinks.Range("A2").Value = 0
$wsPupilBlinks.Range("B2").Value = "Dilley et al."
$wsPupilBlinks.Range("C2").Value = 2020
$wsPupilBlinks.Range("D2").Value = "Visual behaviour in robotic surgery—Demonstrating the validity of the simulated environment"
$wsPupilBlinks.Range("E2").Value = "International Journal of Medical Robotics and Computer Assisted Surgery"
$wsPupilBlinks.Range("F2").Value = "Robotic surgery"
$wsPupilBlinks.Range("G2").Value = "Fundamentals of Robotic Surgery, simulator task"
$wsPupilBlinks.Range("H2").Value = "SDs calculated from inter-quartile ranges (SD = (3/4)*IQR). The paper reports medians."
$wsPupilBlinks.Range("I2").Value = 18
$wsPupilBlinks.Range("J2").Value = 22.7
$wsPupilBlinks.Range("K2").Formula = "=(3/4)*20.87"
$wsPupilBlinks.Range("L2").Value = 14
$wsPupilBlinks.Range("M2").Value = 25.28
$wsPupilBlinks.Range("N2").Formula = "=(3/4)*20.08"
$wsPupilBlinks.Range("O2").Formula = "=SQRT(((I2-1)*POWER(K2,2) + (L2-1)*POWER(N2,2))/((I2-1)+(L2-1)))"
$wsPupilBlinks.Range("P2").Formula = "=(J2-M2)/O2"
$wsPupilBlinks.Range("Q2").Formula = "=P2*(1- (3/(4*(I2+L2)-9)))"
$wsPupilBlinks.Range("R2").Formula = "=SQRT((I2+L2)/(I2*L2)+(POWER(P2,2)/(2*(I2+L2))))"

# --- view/selection state ---
$wsTaskTime.Range("N18").Select() | Out-Null
$wsBimanual.Range("O2:R2").Select() | Out-Null
$wsPupilBlinks.Range("D6").Select() | Out-Null
$wsPupilDilation.Activate() | Out-Null
$wsPupilDilation.Range("R11").Select() | Out-Null

